$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.679.58"
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").Value = "3.096.24"
$ws.Range("E3").Value = "  +1.20%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.00%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "3.097.60"
$ws.Range("E8").Value = "  +1.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.436"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.95%  "

$ws.Range("E10").Value = "  +0.29%  "

$ws.Range("E11").Value = "  +0.79%  "

$ws.Range("E12").Value = "  +3.26%  "

$ws.Range("D13").Value = "3.634.92"
$ws.Range("E13").Value = "  +1.23%  "

$ws.Range("E14").Value = "  +1.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.28%  "

$ws.Range("E16").Value = "  +1.28%  "

$ws.Range("D17").Value = "57.778.98"
$ws.Range("E17").Value = "  +0.91%  "

$ws.Range("D18").Value = "3.103.72"
$ws.Range("E18").Value = "  +1.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.06"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "336.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("E24").Value = "  +1.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.05%  "

$ws.Range("E26").Value = "  -1.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("D28").Value = "0.0₃0908"
$ws.Range("E28").Value = "  +0.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.89%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.46%  "

$ws.Range("E32").Value = "  +2.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.55%  "

$ws.Range("E34").Value = "  +3.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "153.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.75%  "

$ws.Range("E36").Value = "  +3.52%  "

$ws.Range("E37").Value = "  +3.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("E39").Value = "  +1.74%  "

$ws.Range("E40").Value = "  -0.94%  "

$ws.Range("D41").Value = "3.145.08"
$ws.Range("E41").Value = "  +1.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.679"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.52%  "

$ws.Range("E43").Value = "  +0.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "36.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.24%  "

$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("E46").Value = "  +6.93%  "

$ws.Range("D47").Value = "2.278.66"
$ws.Range("E47").Value = "  +0.72%  "

$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.958"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.48%  "
